# Consolidate the "Absent" (column H) report based on the "Real" (column E) column.
# A student is considered Absent (H = 1) for a given date row when the
# "Real" attendance count (column E) is 0; otherwise they are marked
# present (H = 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows run from row 3 to row 21 (row 1 = headers, row 2 = roll/name info).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 3) { $lastRow = 21 }

for ($r = 3; $r -le $lastRow; $r++) {
    $real = $ws.Cells.Item($r, 5).Value2  # Column E = Real
    if ($real -eq 0) {
        $ws.Cells.Item($r, 8).Value = 1  # Column H = Absent
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
